$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.065.23'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.651.58'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.22'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5260'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2595'
$ws.Range('E8').Value = '  -1.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06321'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.34'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07797'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.502'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.650.13'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5486'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0₅8203'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.48'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.075.78'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.568'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.66'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.06'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.024'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.56'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1236'
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.222'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.03'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.430'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05824'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.272'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.547'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.263'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.583'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9473'
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.779'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.409'
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5738'
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01609'
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8434'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.742'
$ws.Range('E40').Value = '  -4.97%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.16'
$ws.Range('E41').Value = '  +3.51%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.029.72'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.794.62'
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '56.96'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4323'
$ws.Range('E47').Value = '  +2.55%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.853'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05143'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.460'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09617'
$ws.Range('E51').Value = '  -0.94%  '
